$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "22÷7=3, 1"
$t.Cell(1, 2).Range.Text = "62÷2=31, 0"
$t.Cell(1, 3).Range.Text = "28÷9=3, 1"
$t.Cell(1, 4).Range.Text = "73÷9=8, 1"
$t.Cell(1, 5).Range.Text = "40÷3=13, 1"

$t.Cell(5, 1).Range.Text = "30÷4=7, 2"
$t.Cell(5, 2).Range.Text = "57÷2=28, 1"
$t.Cell(5, 3).Range.Text = "95÷5=19, 0"
$t.Cell(5, 4).Range.Text = "12÷6=2, 0"
$t.Cell(5, 5).Range.Text = "77÷4=19, 1"

$t.Cell(9, 1).Range.Text = "36÷2=18, 0"
$t.Cell(9, 2).Range.Text = "63÷7=9, 0"
$t.Cell(9, 3).Range.Text = "73÷9=8, 1"
$t.Cell(9, 4).Range.Text = "60÷5=12, 0"
$t.Cell(9, 5).Range.Text = "96÷4=24, 0"

$t.Cell(13, 1).Range.Text = "63÷9=7, 0"
$t.Cell(13, 2).Range.Text = "83÷2=41, 1"
$t.Cell(13, 3).Range.Text = "82÷9=9, 1"
$t.Cell(13, 4).Range.Text = "58÷4=14, 2"
$t.Cell(13, 5).Range.Text = "77÷3=25, 2"

$t.Cell(17, 1).Range.Text = "87÷2=43, 1"
$t.Cell(17, 2).Range.Text = "40÷4=10, 0"
$t.Cell(17, 3).Range.Text = "45÷9=5, 0"
$t.Cell(17, 4).Range.Text = "29÷7=4, 1"
$t.Cell(17, 5).Range.Text = "36÷6=6, 0"
